$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- P4: new column header "2022", formatted like the other year header cells ---
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("P4").Value = 2022

# --- P5: total row, copy format from O5 then switch to new "#,##0.0" number format ---
$ws.Range("O5").Copy()
$ws.Range("P5").PasteSpecial(-4122)
$ws.Range("P5").VerticalAlignment = -4107   # xlBottom (drop explicit vertical centering)
$ws.Range("P5").NumberFormat = "#,##0.0"
$ws.Range("P5").Value = 1188.7

# --- P6: first data row using the "plain" numeric style, creates the 2nd new style ---
$ws.Range("O9").Copy()
$ws.Range("P6").PasteSpecial(-4122)
$ws.Range("P6").VerticalAlignment = -4107
$ws.Range("P6").NumberFormat = "#,##0.0"
$ws.Range("P6").Value = 263.89999999999998

# --- P16: bottom bordered row, creates the 3rd new style ---
$ws.Range("O16").Copy()
$ws.Range("P16").PasteSpecial(-4122)
$ws.Range("P16").VerticalAlignment = -4107
$ws.Range("P16").NumberFormat = "#,##0.0"
$ws.Range("P16").Value = 3.1

# --- P7..P15: reuse the style that was just created for P6 ---
$ws.Range("P6").Copy()
$ws.Range("P7:P15").PasteSpecial(-4122)

$ws.Range("P7").Value = 263.2
$ws.Range("P8").Value = 12.4
$ws.Range("P9").Value = "-"
$ws.Range("P10").Value = 93
$ws.Range("P11").Value = 171.5
$ws.Range("P12").Value = 220.6
$ws.Range("P13").Value = 159.30000000000001
$ws.Range("P14").Value = 1.7
$ws.Range("P15").Value = "-"

$excel.CutCopyMode = 0

# --- leave the active selection on Q7, matching the saved workbook state ---
$ws.Range("Q7").Select()
